$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: capture the current last-row banding style (row 58) and apply it to the new last row (72) ---
$ws.Range("A58:Y58").Copy() | Out-Null
$ws.Range("A72:Y72").PasteSpecial(-4122) | Out-Null

# --- Step 2: re-style row 58 as a normal interior row (it follows the odd/even banding like row 56) ---
$ws.Range("A56:Y56").Copy() | Out-Null
$ws.Range("A58:Y58").PasteSpecial(-4122) | Out-Null

# --- Step 3: style the new interior rows 59-71 by copying the alternating banding from rows 57/56 ---
for ($r = 59; $r -le 71; $r++) {
    if (($r % 2) -eq 1) {
        $ws.Range("A57:Y57").Copy() | Out-Null
    } else {
        $ws.Range("A56:Y56").Copy() | Out-Null
    }
    $ws.Range("A" + $r + ":Y" + $r).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# --- Step 4: write the response values for the new rows (59-72) ---
$ws.Cells.Item(59, 1).Value2 = 45836.60610662037
$rowVals = @('Feminino', '24 a 26 anos', '6º semestre', 'Não', 'Pública', 'Não', 'Até 1 salário mínimo', 'Família', 'Católica', 'Não', 'Sim', 'Não', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não', 'Não', 'Não', 'Não', 'Não', 'Sim')
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(59, 2 + $i).Value2 = $rowVals[$i]
}

$ws.Cells.Item(60, 1).Value2 = 45836.60891416667
$rowVals = @('Masculino', '21 a 23 anos', '8º semestre', 'Sim', 'Pública', 'Não', 'Até 1 salário mínimo', 'Só', 'Católica', 'Não', 'Sim', 'Não', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não', 'Não', 'Não', 'Não', 'Sim', 'Sim')
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(60, 2 + $i).Value2 = $rowVals[$i]
}

$ws.Cells.Item(61, 1).Value2 = 45836.66388969908
$rowVals = @('Masculino', '24 a 26 anos', 'A partir do 10º semestre', 'Sim', 'Privada', 'Não', '4 a 6 salários mínimos', 'Família', 'Evangélica', 'Não', 'Não', 'Não', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não', 'Não', 'Não', 'Sim', 'Não', 'Não')
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(61, 2 + $i).Value2 = $rowVals[$i]
}

$ws.Cells.Item(62, 1).Value2 = 45836.6666487037
$rowVals = @('Prefiro não declarar', '18 a 20 anos', '1º semestre', 'Não', 'Pública', 'Não', 'Até 1 salário mínimo', 'Família', 'Católica', 'Não', 'Sim', 'Não', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não', 'Não', 'Não', 'Não', 'Não', 'Não')
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(62, 2 + $i).Value2 = $rowVals[$i]
}

$ws.Cells.Item(63, 1).Value2 = 45836.66759834491
$rowVals = @('Feminino', '21 a 23 anos', 'A partir do 10º semestre', 'Sim', 'Pública', 'Não', '1 a 2 salários mínimos', 'Família', 'Outros', 'Não', 'Sim', 'Sim', '15 a 17 anos', 'Menos de 1 copo por dia', 'Só em festas', 'Iniciativa própria', 'Liberdade', '“ICE”', 'Não', 'Não', 'Não', 'Sim', 'Sim', 'Não')
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(63, 2 + $i).Value2 = $rowVals[$i]
}

$ws.Cells.Item(64, 1).Value2 = 45836.66823284722
$rowVals = @('Feminino', '21 a 23 anos', '10º semestre', 'Sim', 'Privada', 'Sim', 'Prefiro não declarar', 'Família', 'Outros', 'Não', 'Não', 'Não', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não', 'Não', 'Não', 'Sim', 'Não', 'Sim')
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(64, 2 + $i).Value2 = $rowVals[$i]
}

$ws.Cells.Item(65, 1).Value2 = 45836.6792771412
$rowVals = @('Masculino', '18 a 20 anos', '2º semestre', 'Não', 'Privada', 'Não', 'Prefiro não declarar', 'Família', 'Outros', 'Não', 'Sim', 'Não', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não', 'Não', 'Não', 'Não', 'Sim', 'Sim')
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(65, 2 + $i).Value2 = $rowVals[$i]
}

$ws.Cells.Item(66, 1).Value2 = 45836.6795049537
$rowVals = @('Masculino', '18 a 20 anos', '2º semestre', 'Não', 'Privada', 'Não', '1 a 2 salários mínimos', 'Família', 'Outros', 'Não', 'Sim', 'Não', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Outros', 'Não consumo bebidas alcoólicas', 'Não', 'Não', 'Não', 'Sim', 'Sim', 'Sim')
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(66, 2 + $i).Value2 = $rowVals[$i]
}

$ws.Cells.Item(67, 1).Value2 = 45836.67991284722
$rowVals = @('Masculino', '27 a 29 anos', '9º semestre', 'Sim', 'Pública', 'Sim', '2 a 4 salários mínimos', 'Só', 'Católica', 'Não', 'Sim', 'Sim', '21 a 25 anos', 'Menos de 1 copo por dia', 'Fins de semana', 'Incentivo de amigos', 'Liberdade', 'Cerveja', 'Não', 'Não', 'Não', 'Sim', 'Sim', 'Sim')
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(67, 2 + $i).Value2 = $rowVals[$i]
}

$ws.Cells.Item(68, 1).Value2 = 45836.75218818287
$rowVals = @('Masculino', '18 a 20 anos', '3º semestre', 'Não', 'Pública', 'Não', 'Até 1 salário mínimo', 'Família', 'Católica', 'Não', 'Sim', 'Sim', '15 a 17 anos', '7 ou mais copos por dia', 'Só em festas', 'Iniciativa própria', 'Outros', 'Cerveja', 'Sim', 'Não', 'Sim', 'Sim', 'Sim', 'Sim')
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(68, 2 + $i).Value2 = $rowVals[$i]
}

$ws.Cells.Item(69, 1).Value2 = 45836.753425648145
$rowVals = @('Feminino', '18 a 20 anos', '3º semestre', 'Não', 'Pública', 'Não', 'Prefiro não declarar', 'Amigos', 'Outros', 'Sim', 'Sim', 'Sim', '15 a 17 anos', '3 a 4 copos por dia', 'Só em festas', 'Iniciativa própria', 'Outros', '“ICE”', 'Não', 'Não', 'Não', 'Sim', 'Não', 'Sim')
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(69, 2 + $i).Value2 = $rowVals[$i]
}

$ws.Cells.Item(70, 1).Value2 = 45836.80481577547
$rowVals = @('Prefiro não declarar', '21 a 23 anos', '3º semestre', 'Não', 'Pública', 'Não', 'Até 1 salário mínimo', 'Família', 'Outros', 'Não', 'Não', 'Sim', '21 a 25 anos', 'Não consumo bebidas alcoólicas', 'Outra', 'Iniciativa própria', 'Outros', 'Outros', 'Não', 'Não', 'Não', 'Não', 'Não', 'Não')
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(70, 2 + $i).Value2 = $rowVals[$i]
}

$ws.Cells.Item(71, 1).Value2 = 45836.8887953588
$rowVals = @('Masculino', '24 a 26 anos', '8º semestre', 'Sim', 'Pública', 'Sim', '4 a 6 salários mínimos', 'Família', 'Outros', 'Não', 'Sim', 'Sim', '12 a 14 anos', '5 a 6 copos por dia', 'Fins de semana', 'Incentivo de amigos', 'Outros', 'Cerveja', 'Não', 'Não', 'Sim', 'Sim', 'Sim', 'Não')
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(71, 2 + $i).Value2 = $rowVals[$i]
}

$ws.Cells.Item(72, 1).Value2 = 45836.95655461805
$rowVals = @('Masculino', '21 a 23 anos', '3º semestre', 'Sim', 'Pública', 'Não', 'Prefiro não declarar', 'Família', 'Outros', 'Não', 'Sim', 'Não', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não consumo bebidas alcoólicas', 'Não', 'Não', 'Não', 'Sim', 'Sim', 'Sim')
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(72, 2 + $i).Value2 = $rowVals[$i]
}

# --- Step 5: expand the table (ListObject) to include the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:Y72")) | Out-Null

Write-Host "done"